# Updated cryptos list -- price (D) and volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new text looks like a plain number (e.g. "215.33").
# A bare Value assignment would make Excel auto-convert such literals to a
# numeric cell, but the source data keeps every Price cell as plain text
# (matching unchanged neighbouring cells such as D4 = "1.001"). Format each
# such cell as Text first so the string is preserved, then clear the formatting
# again afterwards so no stray number format lingers behind on the cell.
$numericLookingPriceCells = @("D5", "D6", "D9", "D10", "D11", "D15", "D16", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D48", "D49", "D50", "D51")
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Price (D column) updates
$ws.Range("D2").Value = "25.827.92"
$ws.Range("D3").Value = "1.629.69"
$ws.Range("D13").Value = "1.640.43"
$ws.Range("D14").Value = "1.851.26"
$ws.Range("D17").Value = "0.0₅7517"
$ws.Range("D18").Value = "25.832.31"
$ws.Range("D39").Value = "1.117.32"
$ws.Range("D45").Value = "1.775.19"
$ws.Range("D5").Value = "215.33"
$ws.Range("D6").Value = "0.5112"
$ws.Range("D9").Value = "0.06346"
$ws.Range("D10").Value = "19.51"
$ws.Range("D11").Value = "0.07788"
$ws.Range("D15").Value = "0.5537"
$ws.Range("D16").Value = "63.70"
$ws.Range("D20").Value = "4.437"
$ws.Range("D21").Value = "194.85"
$ws.Range("D22").Value = "9.797"
$ws.Range("D23").Value = "6.019"
$ws.Range("D25").Value = "1.886"
$ws.Range("D26").Value = "141.49"
$ws.Range("D27").Value = "0.1251"
$ws.Range("D28").Value = "15.55"
$ws.Range("D29").Value = "6.727"
$ws.Range("D30").Value = "1.238"
$ws.Range("D31").Value = "0.04872"
$ws.Range("D32").Value = "3.263"
$ws.Range("D33").Value = "3.173"
$ws.Range("D34").Value = "1.546"
$ws.Range("D35").Value = "2.359"
$ws.Range("D36").Value = "0.8972"
$ws.Range("D37").Value = "0.5529"
$ws.Range("D40").Value = "0.01553"
$ws.Range("D41").Value = "1.001"
$ws.Range("D42").Value = "5.558"
$ws.Range("D43").Value = "0.7983"
$ws.Range("D44").Value = "97.45"
$ws.Range("D48").Value = "1.002"
$ws.Range("D49").Value = "54.67"
$ws.Range("D50").Value = "0.05124"
$ws.Range("D51").Value = "7.625"

# Restore default formatting now that the text values are safely stored
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).ClearFormats()
}

# Volume(1h) (E column) updates
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +3.87%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  -7.37%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  +3.68%  "

